# Update cached market-price / profit figures on each Leve sheet.
# Values below were produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -674
$ws.Range("N32").ClearContents()
$ws.Range("H62").Value = 5951.091
$ws.Range("I62").Value = 2938.6
$ws.Range("K62").Value = 2938.6
$ws.Range("M62").Value = -2314.6
$ws.Range("H65").Value = 5951.091
$ws.Range("I65").Value = 2938.6
$ws.Range("K65").Value = 14693
$ws.Range("M65").Value = -11573
$ws.Range("H70").Value = 1193.125
$ws.Range("I70").Value = 681.6667
$ws.Range("K70").Value = 2045.0001
$ws.Range("M70").Value = -1775.0001
$ws.Range("H73").Value = 1193.125
$ws.Range("I73").Value = 681.6667
$ws.Range("K73").Value = 2045.0001
$ws.Range("M73").Value = -1109.0001
$ws.Range("H88").Value = 3999.5
$ws.Range("J88").Value = 3666.3333
$ws.Range("L88").Value = 3666.3333
$ws.Range("N88").Value = -4478.3333
$ws.Range("H91").Value = 3999.5
$ws.Range("J91").Value = 3666.3333
$ws.Range("L91").Value = 3666.3333
$ws.Range("N91").Value = -6474.3333
$ws.Range("H132").Value = 6092.6665
$ws.Range("I132").Value = 6107
$ws.Range("K132").Value = 18321
$ws.Range("M132").Value = -15791

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3999
$ws.Range("J8").Value = 3999
$ws.Range("L8").Value = 3999
$ws.Range("N8").Value = -4287
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H13").Value = 559.75
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 613
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 613
$ws.Range("M13").Value = -256
$ws.Range("N13").Value = -901
$ws.Range("H63").Value = 1927.4286
$ws.Range("I63").Value = 1899.4
$ws.Range("K63").Value = 1899.4
$ws.Range("M63").Value = -1213.4
$ws.Range("H66").Value = 1927.4286
$ws.Range("I66").Value = 1899.4
$ws.Range("K66").Value = 9497
$ws.Range("M66").Value = -6065
$ws.Range("H92").Value = 54966.668
$ws.Range("J92").Value = 54966.668
$ws.Range("L92").Value = 54966.668
$ws.Range("N92").Value = -59958.668

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 91495
$ws.Range("J135").Value = 91495
$ws.Range("L135").Value = 91495
$ws.Range("N135").Value = -101635

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2683.5
$ws.Range("I31").Value = 1863.4
$ws.Range("J31").Value = 4733.75
$ws.Range("K31").Value = 1863.4
$ws.Range("L31").Value = 4733.75
$ws.Range("M31").Value = -1568.4
$ws.Range("N31").Value = -5323.75
$ws.Range("H34").Value = 2683.5
$ws.Range("I34").Value = 1863.4
$ws.Range("J34").Value = 4733.75
$ws.Range("K34").Value = 1863.4
$ws.Range("L34").Value = 4733.75
$ws.Range("M34").Value = -1661.4
$ws.Range("N34").Value = -5137.75
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H106").Value = 9999
$ws.Range("J106").Value = 9999
$ws.Range("L106").Value = 9999
$ws.Range("N106").Value = -12523

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4000000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4000000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4000000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4000224
$ws.Range("H8").Value = 4000000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 4000000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4000000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -4000278
$ws.Range("H11").Value = 3222777.5
$ws.Range("I11").Value = 1000999.4
$ws.Range("K11").Value = 1000999.4
$ws.Range("M11").Value = -1000860.4

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 15000
$ws.Range("I7").Value = 15000
$ws.Range("K7").Value = 15000
$ws.Range("M7").Value = -14888
$ws.Range("H11").Value = 1395
$ws.Range("J11").Value = 1390
$ws.Range("L11").Value = 1390
$ws.Range("N11").Value = -1670
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H29").Value = 15000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H40").Value = 3925.5
$ws.Range("I40").Value = 3925.5
$ws.Range("K40").Value = 3925.5
$ws.Range("M40").Value = -3789.5
$ws.Range("H122").Value = 9698.799999999999
$ws.Range("J122").Value = 9999.5
$ws.Range("L122").Value = 29998.5
$ws.Range("N122").Value = -34898.5
$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 15000
$ws.Range("K126").Value = 45000
$ws.Range("M126").Value = -42530
$ws.Range("H136").Value = 5214.7144
$ws.Range("I136").Value = 5250.5
$ws.Range("K136").Value = 15751.5
$ws.Range("M136").Value = -13201.5
$ws.Range("H139").Value = 45394.5
$ws.Range("J139").Value = 45000
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1226
$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10406
$ws.Range("H63").Value = 28798.8
$ws.Range("J63").Value = 28798.8
$ws.Range("L63").Value = 28798.8
$ws.Range("N63").Value = -30046.8
$ws.Range("H66").Value = 28798.8
$ws.Range("J66").Value = 28798.8
$ws.Range("L66").Value = 86396.39999999999
$ws.Range("N66").Value = -92636.39999999999
